$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F6").Value2 = 56
$ws.Range("F7").Value2 = 1167
$ws.Range("F8").Value2 = 371
$ws.Range("F10").Value2 = 335
$ws.Range("F11").Value2 = 8095
$ws.Range("F13").Value2 = 9657
$ws.Range("F17").Value2 = 475
$ws.Range("F18").Value2 = 6
$ws.Range("F22").Value2 = 282
$ws.Range("F25").Value2 = 49
$ws.Range("F27").Value2 = 385
$ws.Range("F29").Value2 = 1661
$ws.Range("F30").Value2 = 33
$ws.Range("F31").Value2 = 78
$ws.Range("F32").Value2 = 306
$ws.Range("F33").Value2 = 272
$ws.Range("F35").Value2 = 344
$ws.Range("F36").Value2 = 963
$ws.Range("F40").Value2 = 413
$ws.Range("F41").Value2 = 314
$ws.Range("F45").Value2 = 289
$ws.Range("F46").Value2 = 50
$ws.Range("F47").Value2 = 244
$ws.Range("F48").Value2 = 98

$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value2 = 97
$ws.Range("F8").Value2 = 6
$ws.Range("F18").Value2 = 12
$ws.Range("F20").Value2 = 352

$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value2 = 2738

$ws = $wb.Worksheets.Item(4)
$ws.Range("F8").Value2 = 56
$ws.Range("F9").Value2 = 1167
$ws.Range("F10").Value2 = 371
$ws.Range("F14").Value2 = 97
$ws.Range("F16").Value2 = 335
$ws.Range("F17").Value2 = 8095
$ws.Range("F18").Value2 = 9657
$ws.Range("F22").Value2 = 282
$ws.Range("F23").Value2 = 49
$ws.Range("F25").Value2 = 1661
$ws.Range("F26").Value2 = 33
$ws.Range("F27").Value2 = 78
$ws.Range("F28").Value2 = 306
$ws.Range("F29").Value2 = 272
$ws.Range("F31").Value2 = 344
$ws.Range("F33").Value2 = 963
$ws.Range("F40").Value2 = 314
$ws.Range("F42").Value2 = 289
$ws.Range("F43").Value2 = 50
$ws.Range("F44").Value2 = 244
$ws.Range("F45").Value2 = 12
$ws.Range("F47").Value2 = 352
$ws.Range("F48").Value2 = 98
